# Correct problems with fixed interactions in Eduati network
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two fixed interactions currently listed in rows 2 and 3:
# Row 2 was TGFb->TGFRb, Row 3 was EGF->EGFR. Swap the From/To source
# string values so Row 2 becomes EGF->EGFR and Row 3 becomes TGFb->TGFRb.
$ws.Range("A2").Value = "EGF"
$ws.Range("C2").Value = "EGFR"

$ws.Range("A3").Value = "TGFb"
$ws.Range("C3").Value = "TGFRb"

# Update the view: scrolled position and current selection
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("A2:F2").Select()
